# Auto-planned script: content-shuffle across 11 text slots.
# Strategy: Phase 1 replace each slot's current text with a unique
# placeholder marker (so overlapping/cyclic text never collides).
# Phase 2 replace each placeholder with the final destination text.

$d = $word.ActiveDocument
$vt = [char]11

function Replace-Text($old, $new) {
    $range = $d.Content
    $ok = $range.Find.Execute($old, $true, $false, $false, $false, $false, $true, 1, $false, $new, 2)
    if (-not $ok) {
        Write-Output "WARNING: Find failed for: $($old.Substring(0, [Math]::Min(40, $old.Length)))"
    }
    return $ok
}

# ---------- Phase 1: stash current slot contents behind placeholders ----------
Replace-Text ('Apresentar as técnicas experimentais de preparação materialográfica e de caracterização de materiais.') '@@SLOT_A@@' | Out-Null  # was: Objetivos-plain
Replace-Text ('To present the experimental techniques of materialographic preparation and characterization of materials.') '@@SLOT_B@@' | Out-Null  # was: Objetivos-italic
Replace-Text ('6495737 - Durval Rodrigues Junior') '@@SLOT_C@@' | Out-Null  # was: Docente-run0
Replace-Text ('1643715 - Paulo Atsushi Suzuki') '@@SLOT_D@@' | Out-Null  # was: Docente-run1(G block)
Replace-Text ('Difração de raios X. Materialografia. Microscopia óptica. Microscopia eletrônica. Análise térmica.') '@@SLOT_E@@' | Out-Null  # was: ProgResumido-plain
Replace-Text ('X-ray diffraction. Materialography. Optical microscopy. Electron microscopy. Thermal analysis.') '@@SLOT_F@@' | Out-Null  # was: ProgResumido-italic
Replace-Text ('A microestrutura dos materiais. Sistemas e reticulados cristalinos, grupos espaciais e simetria, tipos mais comuns de estruturas cristalinas. Projeção estereográfica. Direção do feixe difratado e a lei de Bragg. Intensidade do feixe difratado. Métodos de difração de raios X. ' + ([char]11) + 'Preparação materialográfica de amostras: corte, embutimento, lixamento e polimento. Técnicas de ataque químico para revelação de fases. Fundamentos de materialografia quantitativa. Microscopia óptica. Técnicas de microscopia eletrônica: varredura e transmissão. Análise química de microrregiões: espectroscopia de energia dispersiva. Técnicas de análise térmica: análise térmica diferencial, calorimetria exploratória diferencial e análise termogravimétrica.') '@@SLOT_G@@' | Out-Null  # was: Programa-plain
Replace-Text ('Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo.') '@@SLOT_H@@' | Out-Null  # was: Avaliacao-Metodo-val
Replace-Text ('Média aritmética de duas provas escritas, testes, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + P2 + TR)/3') '@@SLOT_I@@' | Out-Null  # was: Avaliacao-Criterio-val(K block)
Replace-Text ('Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação') '@@SLOT_J@@' | Out-Null  # was: Avaliacao-Norma-val
Replace-Text ('PADILHA, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985.' + ([char]11) + 'MURPHY, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001.' + ([char]11) + 'WU, Q.; MERCHANT, F.; CASTLEMAN, K. Microscope Image Processing, Academic Press, 2008.' + ([char]11) + 'CULLITY, B. D.; STOCK, S. R. Elements of X-Ray Diffraction, Prentice Hall, 2001.' + ([char]11) + 'YACOBI, B. G.; HOLT, D. B.; KAZMERSKI, L. L. Microanalysis of Solids. Plenum Press, New York, 1994.' + ([char]11) + 'HATAKEYAMA, T.; ZHENHAI, L. Handbook of Thermal Analysis, Wiley, 1999.' + ([char]11) + 'HAINES, P. J. Principles of Thermal Analysis and Calorimetry, Royal Society of Chemistry, 2002.') '@@SLOT_K@@' | Out-Null  # was: Bibliografia-plain

# ---------- Phase 2: drop the correct final text into each placeholder ----------
Replace-Text '@@SLOT_A@@' ('Difração de raios X. Materialografia. Microscopia óptica. Microscopia eletrônica. Análise térmica.') | Out-Null  # now: Objetivos-plain
Replace-Text '@@SLOT_B@@' ('X-ray diffraction. Materialography. Optical microscopy. Electron microscopy. Thermal analysis.') | Out-Null  # now: Objetivos-italic
Replace-Text '@@SLOT_C@@' ('Apresentar as técnicas experimentais de preparação materialográfica e de caracterização de materiais.') | Out-Null  # now: Docente-run0
Replace-Text '@@SLOT_D@@' ('A microestrutura dos materiais. Sistemas e reticulados cristalinos, grupos espaciais e simetria, tipos mais comuns de estruturas cristalinas. Projeção estereográfica. Direção do feixe difratado e a lei de Bragg. Intensidade do feixe difratado. Métodos de difração de raios X. ' + ([char]11) + 'Preparação materialográfica de amostras: corte, embutimento, lixamento e polimento. Técnicas de ataque químico para revelação de fases. Fundamentos de materialografia quantitativa. Microscopia óptica. Técnicas de microscopia eletrônica: varredura e transmissão. Análise química de microrregiões: espectroscopia de energia dispersiva. Técnicas de análise térmica: análise térmica diferencial, calorimetria exploratória diferencial e análise termogravimétrica.') | Out-Null  # now: Docente-run1(G block)
Replace-Text '@@SLOT_E@@' ('Experimentos desenvolvidos em laboratório didático, realização de relatórios para cada experimento e de testes sobre o experimento em estudo.') | Out-Null  # now: ProgResumido-plain
Replace-Text '@@SLOT_F@@' ('To present the experimental techniques of materialographic preparation and characterization of materials.') | Out-Null  # now: ProgResumido-italic
Replace-Text '@@SLOT_G@@' ('Média aritmética de duas provas escritas, testes, trabalhos e relatórios: P1, P2 e TR. Conceito Final = (P1 + P2 + TR)/3') | Out-Null  # now: Programa-plain
Replace-Text '@@SLOT_H@@' ('Aplicação de uma prova escrita e prática dentro do prazo regimental antes do início do próximo semestre letivo. A nota da segunda avaliação será a média aritmética entre a nota da prova de recuperação e a nota final da primeira avaliação') | Out-Null  # now: Avaliacao-Metodo-val
Replace-Text '@@SLOT_I@@' ('PADILHA, A.F. Técnicas de Análise Microestrutural, Ed. Hemus, São Paulo, 1985.' + ([char]11) + 'MURPHY, D. B. Fundamentals of Light Microscopy and Electronic Imaging, Wiley-Liss, 2001.' + ([char]11) + 'WU, Q.; MERCHANT, F.; CASTLEMAN, K. Microscope Image Processing, Academic Press, 2008.' + ([char]11) + 'CULLITY, B. D.; STOCK, S. R. Elements of X-Ray Diffraction, Prentice Hall, 2001.' + ([char]11) + 'YACOBI, B. G.; HOLT, D. B.; KAZMERSKI, L. L. Microanalysis of Solids. Plenum Press, New York, 1994.' + ([char]11) + 'HATAKEYAMA, T.; ZHENHAI, L. Handbook of Thermal Analysis, Wiley, 1999.' + ([char]11) + 'HAINES, P. J. Principles of Thermal Analysis and Calorimetry, Royal Society of Chemistry, 2002.' + ([char]11)) | Out-Null  # now: Avaliacao-Criterio-val(K block)
Replace-Text '@@SLOT_J@@' ('6495737 - Durval Rodrigues Junior') | Out-Null  # now: Avaliacao-Norma-val
Replace-Text '@@SLOT_K@@' ('1643715 - Paulo Atsushi Suzuki') | Out-Null  # now: Bibliografia-plain

Write-Output "done"
